$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new glossary terms at the bottom of the list (rows 19-20).
# Shared strings are appended to the table in the order they are first
# used, so write A20 ("JWT") before A19 ("TOKEN") to reproduce the
# target shared-string ordering (JWT=25, TOKEN=26).
$ws.Range("A20").Value = "JWT"
$ws.Range("A19").Value = "TOKEN"
